$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 9 (Ano 2025) figures per the latest faturamento data
$ws.Range("B9").Value = 3731465.66
$ws.Range("C9").Value = 585140.24
$ws.Range("D9").Value = 4316605.9
$ws.Range("E9").Value = 13.55556317985851
$ws.Range("F9").Value = 86.44443682014149
$ws.Range("G9").Value = -43.44911344774486
$ws.Range("H9").Value = -32.61488796005552
$ws.Range("I9").Value = 37568
$ws.Range("J9").Value = 1595
$ws.Range("K9").Value = 39163
$ws.Range("L9").Value = 27067
$ws.Range("M9").Value = 159.4785495252522
$ws.Range("N9").Value = 8.879253976833045
